$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1764
$ws.Range("F3").Value = 767
$ws.Range("F6").Value = 457
$ws.Range("F7").Value = 1088
$ws.Range("F8").Value = 311
$ws.Range("F10").Value = 100
$ws.Range("F11").Value = 92
$ws.Range("F12").Value = 1075
$ws.Range("F13").Value = 49
$ws.Range("F15").Value = 707
$ws.Range("F16").Value = 776
$ws.Range("F17").Value = 165
$ws.Range("F19").Value = 47
$ws.Range("F20").Value = 605
$ws.Range("F21").Value = 92
$ws.Range("F22").Value = 1683
$ws.Range("F23").Value = 1840
$ws.Range("F24").Value = 457
$ws.Range("F25").Value = 53
$ws.Range("F26").Value = 1686
$ws.Range("F27").Value = 250
$ws.Range("F28").Value = 2501
$ws.Range("F29").Value = 440
$ws.Range("F31").Value = 646
$ws.Range("F32").Value = 125
$ws.Range("F33").Value = 82
$ws.Range("F35").Value = 857
$ws.Range("F36").Value = 1563
$ws.Range("F37").Value = 258
$ws.Range("F39").Value = 504
$ws.Range("F40").Value = 98
$ws.Range("F41").Value = 92
$ws.Range("F42").Value = 132
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 117
$ws.Range("F12").Value = 58
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1764
$ws.Range("F5").Value = 767
$ws.Range("F8").Value = 457
$ws.Range("F9").Value = 1088
$ws.Range("F10").Value = 311
$ws.Range("F12").Value = 100
$ws.Range("F13").Value = 92
$ws.Range("F14").Value = 1075
$ws.Range("F15").Value = 49
$ws.Range("F16").Value = 707
$ws.Range("F17").Value = 776
$ws.Range("F18").Value = 165
$ws.Range("F19").Value = 117
$ws.Range("F20").Value = 117
$ws.Range("F24").Value = 47
$ws.Range("F25").Value = 605
$ws.Range("F26").Value = 92
$ws.Range("F27").Value = 1683
$ws.Range("F28").Value = 1840
$ws.Range("F29").Value = 458
$ws.Range("F30").Value = 53
$ws.Range("F32").Value = 2501
$ws.Range("F33").Value = 440
$ws.Range("F37").Value = 58
$ws.Range("F38").Value = 646
$ws.Range("F39").Value = 125
$ws.Range("F40").Value = 82
$ws.Range("F42").Value = 857
$ws.Range("F43").Value = 1563
$ws.Range("F45").Value = 258
$ws.Range("F46").Value = 504
$ws.Range("F47").Value = 98
$ws.Range("F48").Value = 92
$ws.Range("F49").Value = 132
